$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as text in the source data (e.g. "28.463.92" uses
# dots as thousands separators), so force text format before assigning to avoid
# Excel auto-converting them to numbers/dates.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.463.92'
$ws.Range('E2').Value = '  -0.20%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.566.64'
$ws.Range('E3').Value = '  -2.19%  '

$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.96'
$ws.Range('E5').Value = '  -1.33%  '

$ws.Range('E6').Value = '  -1.28%  '

$ws.Range('E7').Value = '  +0.19%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '46.04'
$ws.Range('E8').Value = '  +4.42%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '24.00'
$ws.Range('E9').Value = '  -0.02%  '

$ws.Range('E10').Value = '  -1.73%  '

$ws.Range('E11').Value = '  -1.55%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0888'
$ws.Range('E12').Value = '  -0.25%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.791.74'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.564.54'
$ws.Range('E14').Value = '  -2.35%  '

$ws.Range('E15').Value = '  -2.77%  '

$ws.Range('E16').Value = '  -2.93%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '28.490.71'
$ws.Range('E17').Value = '  -0.19%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '62.24'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '227.51'
$ws.Range('E19').Value = '  -1.97%  '

$ws.Range('E20').Value = '  -2.67%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0691'
$ws.Range('E21').Value = '  -2.92%  '

$ws.Range('E22').Value = '  +0.12%  '

$ws.Range('E23').Value = '  -5.88%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.13'
$ws.Range('E24').Value = '  -3.21%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.10'
$ws.Range('E25').Value = '  +6.96%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.01'
$ws.Range('E26').Value = '  -1.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.00'
$ws.Range('E27').Value = '  -2.12%  '

$ws.Range('E28').Value = '  -2.67%  '

$ws.Range('E29').Value = '  -4.15%  '

$ws.Range('E30').Value = '  +0.12%  '

$ws.Range('E31').Value = '  -1.65%  '

$ws.Range('E32').Value = '  -4.32%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.21'
$ws.Range('E33').Value = '  -1.49%  '

$ws.Range('E34').Value = '  -3.16%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.394.66'
$ws.Range('E35').Value = '  -1.83%  '

$ws.Range('E37').Value = '  -3.49%  '

$ws.Range('E38').Value = '  +1.47%  '

$ws.Range('E39').Value = '  +2.84%  '

$ws.Range('E40').Value = '  -0.87%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.535'
$ws.Range('E41').Value = '  -1.96%  '

$ws.Range('E42').Value = '  +0.17%  '

$ws.Range('E43').Value = '  -4.45%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.87'
$ws.Range('E44').Value = '  +1.73%  '

$ws.Range('E45').Value = '  -4.29%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.973'
$ws.Range('E46').Value = '  -1.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '62.74'
$ws.Range('E47').Value = '  -3.41%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.704.35'
$ws.Range('E48').Value = '  -2.16%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.86'
$ws.Range('E49').Value = '  -1.91%  '

$ws.Range('E50').Value = '  -0.31%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0103'
$ws.Range('E51').Value = '  -4.54%  '
